# Study 2 and 4: add electricity production share rows to the "z" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("z")

$data = @(
    @("EU27", "Activity", "Production of electricity by biomass and waste", "EU27", "Commodity", "Electricity", "Update", 0.02928176795580111),
    @("EU27", "Activity", "Production of electricity by coal", "EU27", "Commodity", "Electricity", "Update", 0),
    @("EU27", "Activity", "Production of electricity by gas", "EU27", "Commodity", "Electricity", "Update", 0.04392265193370166),
    @("EU27", "Activity", "Production of electricity by hydro", "EU27", "Commodity", "Electricity", "Update", 0.05524861878453039),
    @("EU27", "Activity", "Production of electricity by nuclear", "EU27", "Commodity", "Electricity", "Update", 0.1325966850828729),
    @("EU27", "Activity", "Production of electricity by petroleum and other oil derivatives", "EU27", "Commodity", "Electricity", "Update", 0),
    @("EU27", "Activity", "Production of electricity by solar photovoltaic", "EU27", "Commodity", "Electricity", "Update", 0.1850828729281768),
    @("EU27", "Activity", "Production of electricity by wind", "EU27", "Commodity", "Electricity", "Update", 0.5538674033149171)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
